$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

for ($i = 1; $i -le 10; $i++) {
    $row = 11 + $i
    $ws.Cells.Item($row, 2).Value = $i
}
